$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).ColumnWidth = 14.833333333333334
$ws.Columns.Item(2).ColumnWidth = 15.666666666666666

$ws.Cells.Item(1, 1).Value = -0.11720183821370256
$ws.Cells.Item(1, 2).Value = 0.11713823926989164
$ws.Cells.Item(2, 1).Value = -0.13282508401517923
$ws.Cells.Item(2, 2).Value = 0.13254847065938336
$ws.Cells.Item(3, 1).Value = -0.084656897551852239
$ws.Cells.Item(3, 2).Value = 0.084443168323554474
$ws.Cells.Item(4, 1).Value = -0.076443168346566281
$ws.Cells.Item(4, 2).Value = 0.075987621939674455
$ws.Cells.Item(5, 1).Value = -0.072987621952059989
$ws.Cells.Item(5, 2).Value = 0.071440318927064261
$ws.Cells.Item(6, 1).Value = -0.03419113132425089
$ws.Cells.Item(6, 2).Value = 0.033778127179425965
$ws.Cells.Item(7, 1).Value = -0.023778127211692812
$ws.Cells.Item(7, 2).Value = 0.023682446692135262
$ws.Cells.Item(8, 1).Value = -0.013682446725082009
$ws.Cells.Item(8, 2).Value = 0.013525665254010999
$ws.Cells.Item(9, 1).Value = -0.011525665268716345
$ws.Cells.Item(9, 2).Value = 0.011401874279278523
$ws.Cells.Item(10, 1).Value = -0.0094018742945340961
$ws.Cells.Item(10, 2).Value = 0.0093939653421610814
$ws.Cells.Item(11, 1).Value = -0.0063939653597930857
$ws.Cells.Item(11, 2).Value = 0.00638184143305498
$ws.Cells.Item(12, 1).Value = -0.0028818414519626856
$ws.Cells.Item(12, 2).Value = 0.002800990237207035
$ws.Cells.Item(13, 1).Value = 0.00069900974367254065
$ws.Cells.Item(13, 2).Value = -0.00073159615661388955
$ws.Cells.Item(14, 1).Value = 0.0087315961268652487
$ws.Cells.Item(14, 2).Value = -0.0087411182954832967
$ws.Cells.Item(15, 1).Value = 0.0097411182824824749
$ws.Cells.Item(15, 2).Value = -0.0097468508294697997
$ws.Cells.Item(16, 1).Value = 0.011746850814359888
$ws.Cells.Item(16, 2).Value = -0.011763458007545413
$ws.Cells.Item(17, 1).Value = -0.0040036077341936149
$ws.Cells.Item(17, 2).Value = 0.0039999999807287523
$ws.Cells.Item(18, 1).Value = 0.021686844730751886
$ws.Cells.Item(18, 2).Value = -0.021759051558333908
$ws.Cells.Item(19, 1).Value = -0.012091801542066349
$ws.Cells.Item(19, 2).Value = 0.012017096322195275
$ws.Cells.Item(20, 1).Value = -0.0080170963326438027
$ws.Cells.Item(20, 2).Value = 0.0080056971933615273
$ws.Cells.Item(21, 1).Value = -0.0040056972039366201
$ws.Cells.Item(21, 2).Value = 0.0039999999893360894
$ws.Cells.Item(22, 1).Value = -0.043891573119896776
$ws.Cells.Item(22, 2).Value = 0.043680868355872704
$ws.Cells.Item(23, 1).Value = -0.03868086837135376
$ws.Cells.Item(23, 2).Value = 0.038286550203872949
$ws.Cells.Item(24, 1).Value = -0.018286550256284784
$ws.Cells.Item(24, 2).Value = 0.018188749791983483
$ws.Cells.Item(25, 1).Value = -0.035249187615820077
$ws.Cells.Item(25, 2).Value = 0.035210958202268827
$ws.Cells.Item(26, 1).Value = -0.032710958216712882
$ws.Cells.Item(26, 2).Value = 0.032665745961518056
$ws.Cells.Item(27, 1).Value = -0.03016574597624011
$ws.Cells.Item(27, 2).Value = 0.029922478691206233
$ws.Cells.Item(28, 1).Value = -0.027922478706250864
$ws.Cells.Item(28, 2).Value = 0.027774347499216923
$ws.Cells.Item(29, 1).Value = -0.051294208183182199
$ws.Cells.Item(29, 2).Value = 0.051181206741716778
$ws.Cells.Item(30, 1).Value = 0.00881879310410838
$ws.Cells.Item(30, 2).Value = -0.0089414776706084886
$ws.Cells.Item(31, 1).Value = 0.015941477642051893
$ws.Cells.Item(31, 2).Value = -0.015965682403404458
$ws.Cells.Item(32, 1).Value = 0.025965682367933951
$ws.Cells.Item(32, 2).Value = -0.025990885742677605
